# fix(publipostage): Refactor synthetic array /3
#
# The "statut" emoji codes (column A) and the color name that went with
# the black square (column B, "statut_label") are being refreshed from a
# black/red/orange/green square palette to a blue/red/orange/green book
# palette:
#   ⬛ -> 📘   (and its label "noir" -> "bleu")
#   🟥 -> 📕   (label "rouge" is unchanged)
#   🟧 -> 📙   (label "orange" is unchanged)
#   🟩 -> 📗   (label "vert" is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

[void]$ws.Cells.Replace("⬛", "📘", $xlWhole)
[void]$ws.Cells.Replace("🟥", "📕", $xlWhole)
[void]$ws.Cells.Replace("🟧", "📙", $xlWhole)
[void]$ws.Cells.Replace("🟩", "📗", $xlWhole)

[void]$ws.Cells.Replace("noir", "bleu", $xlWhole)
